# noteshare_data_schema.xlsx — "backend dev push 1"
#
# 1. "books" sheet gains a new "page_count" column, inserted right before
#    the existing "created_at" column (so the old column G shifts to H).
# 2. Selection/active-tab state moves: "books" becomes the active sheet
#    (it was "patrons"), with cell F10 selected. "book_likes" also has its
#    selection moved to F10 (without becoming the active tab).

$wb = $excel.ActiveWorkbook

$wsBooks = $wb.Worksheets.Item("books")
$wsBookLikes = $wb.Worksheets.Item("book_likes")

# --- books: insert a new column before "created_at" (currently column G) ---
[void]$wsBooks.Columns.Item(7).Insert()
$wsBooks.Range("G1").Value = "page_count"

# --- book_likes: selection moves to F10, sheet itself is not activated ---
[void]$wsBookLikes.Range("F10").Select()

# --- books: becomes the active sheet, with F10 selected ---
# (selecting a range on a sheet activates that sheet, which is exactly the
#  desired end state: activeTab points at "books", and "patrons" loses its
#  previous tabSelected flag while keeping its own last selection.)
[void]$wsBooks.Range("F10").Select()
